$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 730.913
$ws.Range("J17").Value = 718.6818
$ws.Range("L17").Value = 2156.0454
$ws.Range("N17").Value = -2492.0454
$ws.Range("H113").Value = 5867.9355
$ws.Range("I113").Value = 6045.1875
$ws.Range("J113").Value = 5678.8667
$ws.Range("K113").Value = 6045.1875
$ws.Range("L113").Value = 5678.8667
$ws.Range("M113").Value = -2791.1875
$ws.Range("N113").Value = -12186.8667
$ws.Range("H125").Value = 34332.3
$ws.Range("I125").Value = 104577
$ws.Range("J125").Value = 4227.4287
$ws.Range("K125").Value = 941193
$ws.Range("L125").Value = 38046.85830000001
$ws.Range("M125").Value = -938733
$ws.Range("N125").Value = -42966.85830000001
$ws.Range("H132").Value = 16546.305
$ws.Range("I132").Value = 17087.045
$ws.Range("K132").Value = 51261.13499999999
$ws.Range("M132").Value = -48731.13499999999
$ws.Range("H138").Value = 24263.4
$ws.Range("J138").Value = 86946.836
$ws.Range("L138").Value = 260840.508
$ws.Range("N138").Value = -271120.508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3522.611
$ws.Range("I45").Value = 1770.8182
$ws.Range("J45").Value = 6275.4287
$ws.Range("K45").Value = 1770.8182
$ws.Range("L45").Value = 6275.4287
$ws.Range("M45").Value = -1393.8182
$ws.Range("N45").Value = -7029.4287
$ws.Range("H97").Value = 1250.258
$ws.Range("I97").Value = 1017.12
$ws.Range("K97").Value = 1017.12
$ws.Range("M97").Value = -521.12
$ws.Range("H110").Value = 5714.07
$ws.Range("I110").Value = 5728.75
$ws.Range("J110").Value = 5518.3335
$ws.Range("K110").Value = 5728.75
$ws.Range("L110").Value = 5518.3335
$ws.Range("M110").Value = -3683.75
$ws.Range("N110").Value = -9608.333500000001
$ws.Range("H132").Value = 1410.8667
$ws.Range("I132").Value = 1012.53845
$ws.Range("K132").Value = 3037.61535
$ws.Range("M132").Value = -507.61535
$ws.Range("H135").Value = 49475
$ws.Range("J135").Value = 49475
$ws.Range("L135").Value = 49475
$ws.Range("N135").Value = -59615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 708.6667
$ws.Range("I80").Value = 767
$ws.Range("K80").Value = 767
$ws.Range("M80").Value = 231
$ws.Range("H83").Value = 708.6667
$ws.Range("I83").Value = 767
$ws.Range("K83").Value = 3835
$ws.Range("M83").Value = 1157

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3847190.2
$ws.Range("I31").Value = 4167498
$ws.Range("J31").Value = 3499.5
$ws.Range("K31").Value = 4167498
$ws.Range("L31").Value = 3499.5
$ws.Range("M31").Value = -4167203
$ws.Range("N31").Value = -4089.5
$ws.Range("H34").Value = 3847190.2
$ws.Range("I34").Value = 4167498
$ws.Range("J34").Value = 3499.5
$ws.Range("K34").Value = 4167498
$ws.Range("L34").Value = 3499.5
$ws.Range("M34").Value = -4167296
$ws.Range("N34").Value = -3903.5
$ws.Range("H99").Value = 4332.3335
$ws.Range("I99").Value = 3771.8
$ws.Range("J99").Value = 5733.6665
$ws.Range("K99").Value = 3771.8
$ws.Range("L99").Value = 5733.6665
$ws.Range("M99").Value = -2273.8
$ws.Range("N99").Value = -8729.666499999999
$ws.Range("H126").Value = 4332.3335
$ws.Range("I126").Value = 3771.8
$ws.Range("J126").Value = 5733.6665
$ws.Range("K126").Value = 11315.4
$ws.Range("L126").Value = 17200.9995
$ws.Range("M126").Value = -8845.400000000001
$ws.Range("N126").Value = -22140.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 91578296
$ws.Range("I4").Value = 96044610
$ws.Range("J4").Value = 84878824
$ws.Range("K4").Value = 288133830
$ws.Range("L4").Value = 254636472
$ws.Range("M4").Value = -288133718
$ws.Range("N4").Value = -254636696
$ws.Range("H7").Value = 296.25
$ws.Range("I7").Value = 295
$ws.Range("K7").Value = 885
$ws.Range("M7").Value = -773
$ws.Range("H12").Value = 246.38889
$ws.Range("J12").Value = 267.83334
$ws.Range("L12").Value = 803.5000200000001
$ws.Range("N12").Value = -1149.50002
$ws.Range("H58").Value = 4999.476
$ws.Range("I58").Value = 4994
$ws.Range("J58").Value = 4999.75
$ws.Range("K58").Value = 14982
$ws.Range("L58").Value = 14999.25
$ws.Range("M58").Value = -14854
$ws.Range("N58").Value = -15255.25
$ws.Range("H126").Value = 2014.5
$ws.Range("I126").Value = 2014.5
$ws.Range("K126").Value = 6043.5
$ws.Range("M126").Value = -1103.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2541.8572
$ws.Range("I126").Value = 1965.5
$ws.Range("K126").Value = 5896.5
$ws.Range("M126").Value = -3426.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1691.0555
$ws.Range("I16").Value = 1277.5
$ws.Range("J16").Value = 4999.5
$ws.Range("K16").Value = 1277.5
$ws.Range("L16").Value = 4999.5
$ws.Range("M16").Value = -1107.5
$ws.Range("N16").Value = -5339.5
$ws.Range("H46").Value = 2092.5625
$ws.Range("J46").Value = 2634.6843
$ws.Range("L46").Value = 2634.6843
$ws.Range("N46").Value = -3010.6843
$ws.Range("H63").Value = 74999
$ws.Range("J63").Value = 74999
$ws.Range("L63").Value = 74999
$ws.Range("N63").Value = -76497
$ws.Range("H66").Value = 74999
$ws.Range("J66").Value = 74999
$ws.Range("L66").Value = 224997
$ws.Range("N66").Value = -232485
$ws.Range("H75").Value = 150000
$ws.Range("J75").Value = 150000
$ws.Range("L75").Value = 150000
$ws.Range("N75").Value = -151872
$ws.Range("H78").Value = 150000
$ws.Range("J78").Value = 150000
$ws.Range("L78").Value = 450000
$ws.Range("N78").Value = -459360
$ws.Range("H88").Value = 64499.6
$ws.Range("I88").Value = 54999
$ws.Range("J88").Value = 66874.75
$ws.Range("K88").Value = 54999
$ws.Range("L88").Value = 66874.75
$ws.Range("M88").Value = -54571
$ws.Range("N88").Value = -67730.75
$ws.Range("H91").Value = 64499.6
$ws.Range("I91").Value = 54999
$ws.Range("J91").Value = 66874.75
$ws.Range("K91").Value = 54999
$ws.Range("L91").Value = 66874.75
$ws.Range("M91").Value = -53517
$ws.Range("N91").Value = -69838.75
$ws.Range("H122").Value = 3229.8462
$ws.Range("I122").Value = 3600
$ws.Range("J122").Value = 2998.5
$ws.Range("K122").Value = 10800
$ws.Range("L122").Value = 8995.5
$ws.Range("M122").Value = -8350
$ws.Range("N122").Value = -13895.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9643389
$ws.Range("I122").Value = 11396223
$ws.Range("J122").Value = 2801
$ws.Range("K122").Value = 34188669
$ws.Range("L122").Value = 8403
$ws.Range("M122").Value = -34186219
$ws.Range("N122").Value = -13303
$ws.Range("H126").Value = 628979.75
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 1004267.6
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 3012802.8
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -3017742.8
$ws.Range("H136").Value = 16723.373
$ws.Range("I136").Value = 20761
$ws.Range("K136").Value = 62283
$ws.Range("M136").Value = -59733
